# Swap "System, <email>" to "<email>, System" in the "Recorded By" column (G)
# wherever it currently reads exactly "System, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

$rng = $ws.Range("G1:G" + $lastRow)
$rng.Replace($oldValue, $newValue)
